$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 7: V = 28 - 1.55
$ws.Range("A7").Value = "V"
$ws.Range("B7").Formula = "=28-1.55"

# New row 8: I = 0.01 A
$ws.Range("A8").Value = "I"
$ws.Range("B8").Value = 0.01
$ws.Range("C8").Value = "A"

# New row 9: R = B7/B8, Ohm
$ws.Range("A9").Value = "R"
$ws.Range("B9").Formula = "=B7/B8"
$ws.Range("C9").Value = "Ohm"

# Update selection to match target state
$ws.Range("B2").Select()
